# The edit renames the internal drawing-object names (wp:docPr/@name and
# pic:cNvPr/@name) of the four logo pictures embedded in the document's
# headers/footers:
#   - the Pearson logo pictures (in footer1.xml / footer2.xml) go from
#     "image2.png" -> "image1.png"
#   - the BTEC logo pictures (in header1.xml / header2.xml) go from
#     "image1.jpg" -> "image2.jpg"
#
# InlineShape has no writable Name property in the Word object model, so
# the rename is applied through the document's WordOpenXML round-trip,
# which exposes the underlying header/footer part markup (including the
# drawing name attributes) as editable text.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image2.png"', 'name="image1.png"')
$xml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')

$d.WordOpenXML = $xml
